$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.350.16'
$ws.Range('E2').Value = '  -2.74%  '
$ws.Range('D3').Value = '1.646.10'
$ws.Range('E3').Value = '  -4.30%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'0.9998"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.51%  '
$ws.Range('D6').Value = "'306.00"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D7').Value = "'0.3601"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.66%  '
$ws.Range('D8').Value = "'47.59"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.92%  '
$ws.Range('D9').Value = "'0.3270"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.65%  '
$ws.Range('D10').Value = "'1.115"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.88%  '
$ws.Range('D11').Value = "'0.06895"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -7.90%  '
$ws.Range('D12').Value = "'0.9990"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = "'5.916"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.72%  '
$ws.Range('D14').Value = "'19.12"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -8.86%  '
$ws.Range('D15').Value = '1.649.26'
$ws.Range('E15').Value = '  -4.15%  '
$ws.Range('D16').Value = "'6.567"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -6.25%  '
$ws.Range('D17').Value = "'0.00001036"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.07%  '
$ws.Range('D18').Value = "'0.06485"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.82%  '
$ws.Range('D19').Value = "'0.9999"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.38%  '
$ws.Range('D20').Value = "'76.29"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -9.55%  '
$ws.Range('D21').Value = "'5.896"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.86%  '
$ws.Range('D22').Value = "'15.60"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -10.16%  '
$ws.Range('D23').Value = "'12.16"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -9.33%  '
$ws.Range('D24').Value = '24.358.41'
$ws.Range('E24').Value = '  -2.94%  '
$ws.Range('D25').Value = "'2.418"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.54%  '
$ws.Range('D26').Value = "'2.289"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -18.99%  '
$ws.Range('D27').Value = "'145.65"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.24%  '
$ws.Range('D28').Value = "'18.20"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -11.46%  '
$ws.Range('D29').Value = '1.827.47'
$ws.Range('E29').Value = '  -4.38%  '
$ws.Range('D30').Value = "'124.09"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.26%  '
$ws.Range('D31').Value = "'1.146"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.37%  '
$ws.Range('D32').Value = "'4.039"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('D33').Value = "'5.535"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -19.63%  '
$ws.Range('B34').Value = 'WEMIXTOKEN'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').Value = "'1.678"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.26%  '
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').Value = "'0.08308"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.03%  '
$ws.Range('D36').Value = "'12.21"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -11.69%  '
$ws.Range('D37').Value = "'5.176"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.51%  '
$ws.Range('D38').Value = "'0.05995"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.03%  '
$ws.Range('D39').Value = "'0.02203"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.34%  '
$ws.Range('D40').Value = "'1.200"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.07%  '
$ws.Range('D41').Value = "'8.166"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -9.00%  '
$ws.Range('D42').Value = "'0.2032"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.20%  '
$ws.Range('D43').Value = "'0.9995"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('D44').Value = "'0.5774"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -10.69%  '
$ws.Range('D45').Value = "'3.718"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.88%  '
$ws.Range('D46').Value = "'12.58"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -10.21%  '
$ws.Range('D47').Value = "'0.5516"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -10.47%  '
$ws.Range('D48').Value = "'121.41"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.32%  '
$ws.Range('D49').Value = "'1.923"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.57%  '
$ws.Range('D50').Value = "'0.06879"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.12%  '
$ws.Range('D51').Value = "'73.82"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.81%  '
